$d = $word.ActiveDocument

# Sprint No. table cell value: "1" -> "2"
$d.Content.Find.Execute("1", $true, $true, $false, $false, $false, $true, 1, $false, "2", 1)

# Review Date table cell value: "02/09/18" -> "02/21/18"
$d.Content.Find.Execute("02/09/18", $true, $true, $false, $false, $false, $true, 1, $false, "02/21/18", 1)
